# "one more month down" -- roll the workbook forward so the September
# prediction sheet is populated with this month's tallies and October
# becomes the active sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Day of Week Pred Sept (sheet index 8): fill in the day-of-week
#    transaction counts that were collected for the month, and clear
#    out the now-unused day-31 placeholder row (September has 30 days).
# ---------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item(8)

$ws8.Range("H2").Value = 5
$ws8.Range("B3").Value = 9
$ws8.Range("C4").Value = 8
$ws8.Range("D5").Value = 7
$ws8.Range("E6").Value = 4
$ws8.Range("F7").Value = 8
$ws8.Range("G8").Value = 11
$ws8.Range("H9").Value = 8
$ws8.Range("B10").Value = 11
$ws8.Range("C11").Value = 6
$ws8.Range("D12").Value = 6
$ws8.Range("E13").Value = 4
$ws8.Range("F14").Value = 5
$ws8.Range("G15").Value = 3
$ws8.Range("H16").Value = 5
$ws8.Range("B17").Value = 7
$ws8.Range("C18").Value = 11
$ws8.Range("D19").Value = 6
$ws8.Range("E20").Value = 4
$ws8.Range("F21").Value = 8
$ws8.Range("G22").Value = 4
$ws8.Range("H23").Value = 8
$ws8.Range("B24").Value = 7
$ws8.Range("C25").Value = 11
$ws8.Range("D26").Value = 4
$ws8.Range("E27").Value = 4
$ws8.Range("F28").Value = 11
$ws8.Range("G29").Value = 3
$ws8.Range("H30").Value = 8
$ws8.Range("B31").Value = 10

# Row 32 (the "31st day" placeholder, A32=31) no longer applies -- remove it.
$ws8.Range("A32").ClearContents()

# ---------------------------------------------------------------------
# 2. Tidy up selections left over on the other monthly sheets.
#    (Activating a sheet while selecting a range updates that sheet's
#    stored cursor position without disturbing its data.)
# ---------------------------------------------------------------------

# March sheet: selection moves from B1 to C1.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C1").Select() | Out-Null

# August sheet: selection moves from G11 to G1 (no longer the active tab).
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("G1").Select() | Out-Null

# September sheet: leave the cursor on the new last row/entry.
$ws8.Range("A32").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. October becomes the newly active sheet/tab.
# ---------------------------------------------------------------------
$ws9 = $wb.Worksheets.Item(9)
$ws9.Activate() | Out-Null
$ws9.Range("A6").Select() | Out-Null
